$wb = $excel.ActiveWorkbook
try {
    Write-Host ("Win TabRatio=" + $excel.ActiveWindow.TabRatio)
} catch { Write-Host ("ERR win " + $_) }
try {
    Write-Host ("App TabRatio=" + $excel.TabRatio)
} catch { Write-Host ("ERR app " + $_) }
try {
    $excel.ActiveWindow.TabRatio = 0.5
    Write-Host ("Win TabRatio after=" + $excel.ActiveWindow.TabRatio)
} catch { Write-Host ("ERR win set " + $_) }
